$p = $ppt.ActivePresentation

# --- Update the "datetimeFigureOut" date placeholders from 6/22/2019 to 6/23/2019 ---
# Slide Master
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "6/22/2019") {
        $shp.TextFrame.TextRange.Text = "6/23/2019"
    }
}

# Slide Layouts
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "6/22/2019") {
            $shp.TextFrame.TextRange.Text = "6/23/2019"
        }
    }
}

# Notes Master
$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $shp = $nm.Shapes.Item($i)
    if ($shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "6/22/2019") {
        $shp.TextFrame.TextRange.Text = "6/23/2019"
    }
}

# --- Slide 13: Table 11, cell (row1,col1) text "- 6.77" -> "- .`677" ---
$s13 = $p.Slides.Item(13)
$tbl13 = $s13.Shapes.Item(8).Table
$tbl13.Cell(1,1).Shape.TextFrame.TextRange.Text = "- .``677"

# --- Slide 14: TextBox 2, paragraph run " (1.2480) <  " -> " (1.2480) >  " ---
$s14 = $p.Slides.Item(14)
$shp14 = $s14.Shapes.Item(2)
$tr14 = $shp14.TextFrame.TextRange
$found = $tr14.Replace(" (1.2480) <  ", " (1.2480) >  ")
